# Updated symbol list on Mon Dec 19 18:45:17 UTC 2022 with GitHub Actions
# Refresh the "Price" (column D) quotes and two "Volume(1h)" (column E)
# labels to match the latest pull of the coinranking symbol list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices stored as plain text (not numbers), so every
# numeric-looking replacement is written with a leading apostrophe to keep
# it a text value instead of letting Excel coerce it into a Number cell.
$priceUpdates = @{
    "D2"  = "242.41"
    "D3"  = "21.52"
    "D4"  = "5.316"
    "D5"  = "0.05631"
    "D6"  = "3.377"
    "D7"  = "6.377"
    "D8"  = "0.8067"
    "D9"  = "0.9575"
    "D10" = "0.1421"
    "D11" = "0.07484"
    "D12" = "0.03213"
    "D13" = "0.03069"
    "D14" = "0.09279"
    "D15" = "3.572"
    "D16" = "0.001645"
    "D17" = "0.04710"
    "D18" = "0.0005825"
    "D19" = "0.006357"
    "D20" = "0.004982"
    "D22" = "0.0001501"
    "D23" = "0.0003103"
    "D24" = "3.769"
    "D25" = "2.098"
    "D26" = "0.3254"
    "D40" = "0.03933"
    "D41" = "0.006979"
    "D42" = "0.003503"
    "D43" = "0.1033"
    "D44" = "0.007482"
    "D45" = "0.00005944"
    "D46" = "0.00000000751"
    "D47" = "0.0005505"
    "D48" = "0.6830"
    "D49" = "0.05167"
    "D50" = "0.00002102"
    "D51" = "0.01011"
}

foreach ($addr in $priceUpdates.Keys) {
    $ws.Range($addr).Value = "'" + $priceUpdates[$addr]
}

# Column E "Worstin24h" / "Bestin24h" badges moved to different rows.
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E47").Value = "46ACDXExchangeACXTWorstin24h"
